$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 2669
$ws.Range('J3').Value = 2747
$ws.Range('B4').Value = 1674
$ws.Range('C4').Value = 1822
$ws.Range('I4').Value = 1759
$ws.Range('J4').Value = 620
$ws.Range('J5').Value = 214
$ws.Range('J6').Value = 3362
$ws.Range('B7').Value = 23306
$ws.Range('C7').Value = 28365
$ws.Range('I7').Value = 26205
$ws.Range('J7').Value = 9612

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J3').Value = 23
$ws.Range('J4').Value = 11
$ws.Range('J7').Value = 101

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J3').Value = 35
$ws.Range('J7').Value = 115

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('J3').Value = 9
$ws.Range('J7').Value = 34

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 44
$ws.Range('J7').Value = 134

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 74
$ws.Range('J6').Value = 104
$ws.Range('J7').Value = 346

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J3').Value = 18
$ws.Range('J7').Value = 73

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 74
$ws.Range('J3').Value = 70
$ws.Range('J6').Value = 90
$ws.Range('J7').Value = 251

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J6').Value = 93
$ws.Range('J7').Value = 289
$ws.Range('J8').Value = 601
$ws.Range('J9').Value = 59
$ws.Range('J10').Value = 55
$ws.Range('J11').Value = 140
$ws.Range('J14').Value = 34
$ws.Range('J15').Value = 115
$ws.Range('J19').Value = 308
$ws.Range('J20').Value = 199
$ws.Range('J22').Value = 23
$ws.Range('J23').Value = 98
$ws.Range('J27').Value = 55
$ws.Range('J29').Value = 557
$ws.Range('J31').Value = 73
$ws.Range('J33').Value = 391
$ws.Range('J34').Value = 51
$ws.Range('J36').Value = 145
$ws.Range('J41').Value = 65
$ws.Range('J42').Value = 373
$ws.Range('J46').Value = 33
$ws.Range('J49').Value = 61
$ws.Range('J52').Value = 242
$ws.Range('J53').Value = 92
$ws.Range('J54').Value = 190
$ws.Range('B63').Value = 378
$ws.Range('C63').Value = 253
$ws.Range('I63').Value = 211
$ws.Range('J63').Value = 46
$ws.Range('J64').Value = 64
$ws.Range('J65').Value = 251
$ws.Range('J67').Value = 346
$ws.Range('J70').Value = 18
$ws.Range('J76').Value = 134
$ws.Range('J77').Value = 82
$ws.Range('J78').Value = 133
$ws.Range('J83').Value = 228
$ws.Range('J85').Value = 445
$ws.Range('J86').Value = 58
$ws.Range('J88').Value = 97
$ws.Range('J89').Value = 101
$ws.Range('J90').Value = 106
$ws.Range('J94').Value = 82
$ws.Range('J95').Value = 151
$ws.Range('J96').Value = 115
$ws.Range('J98').Value = 58
$ws.Range('J99').Value = 134
$ws.Range('B101').Value = 23306
$ws.Range('C101').Value = 28365
$ws.Range('I101').Value = 26205
$ws.Range('J101').Value = 9612

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J2').Value = 68
$ws.Range('J3').Value = 79
$ws.Range('J7').Value = 228

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J2').Value = 56
$ws.Range('J3').Value = 44
$ws.Range('J7').Value = 151

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 104
$ws.Range('J7').Value = 391

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('J6').Value = 29
$ws.Range('J7').Value = 61

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 47
$ws.Range('J3').Value = 39
$ws.Range('J4').Value = 16
$ws.Range('J7').Value = 190

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 163
$ws.Range('J7').Value = 557

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 76
$ws.Range('J7').Value = 308

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J6').Value = 75
$ws.Range('J7').Value = 134

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J3').Value = 170
$ws.Range('J6').Value = 128
$ws.Range('J7').Value = 445

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J3').Value = 27
$ws.Range('J7').Value = 93

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J5').Value = 1
$ws.Range('J7').Value = 65

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J6').Value = 189
$ws.Range('J7').Value = 373

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J3').Value = 14
$ws.Range('J7').Value = 55

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J3').Value = 46
$ws.Range('J7').Value = 133

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('J2').Value = 11
$ws.Range('J7').Value = 33

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('J2').Value = 28
$ws.Range('J7').Value = 98

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('J2').Value = 20
$ws.Range('J7').Value = 64

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 64
$ws.Range('J3').Value = 60
$ws.Range('J6').Value = 49
$ws.Range('J7').Value = 199

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J2').Value = 53
$ws.Range('J3').Value = 38
$ws.Range('J7').Value = 145

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J3').Value = 67
$ws.Range('J4').Value = 12
$ws.Range('J6').Value = 105
$ws.Range('J7').Value = 242

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('J3').Value = 13
$ws.Range('J6').Value = 19
$ws.Range('J7').Value = 51

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J2').Value = 19
$ws.Range('J7').Value = 82

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J3').Value = 31
$ws.Range('J7').Value = 115

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J2').Value = 17
$ws.Range('J7').Value = 58

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 50
$ws.Range('J4').Value = 9
$ws.Range('J6').Value = 52
$ws.Range('J7').Value = 140

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('J3').Value = 16
$ws.Range('J7').Value = 59

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('J2').Value = 9
$ws.Range('J7').Value = 18

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J2').Value = 22
$ws.Range('J6').Value = 42
$ws.Range('J7').Value = 97

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 189
$ws.Range('J3').Value = 194
$ws.Range('J6').Value = 172
$ws.Range('J7').Value = 601

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('J3').Value = 13
$ws.Range('J7').Value = 55

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('J4').Value = 29
$ws.Range('J7').Value = 58

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J2').Value = 36
$ws.Range('J3').Value = 28
$ws.Range('J4').Value = 2
$ws.Range('J7').Value = 106

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J2').Value = 16
$ws.Range('J7').Value = 92

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('J2').Value = 13
$ws.Range('J7').Value = 23

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J2').Value = 26
$ws.Range('J6').Value = 16
$ws.Range('J7').Value = 82

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J3').Value = 83
$ws.Range('J4').Value = 7
$ws.Range('J7').Value = 289
